# translations.xlsx update: add login / OTP / verify-screen translation rows
# (rows 38-80) to Sheet1, matching the "login screen, otp screen, verify screens
# are changed correctly" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 keeps its pre-existing wrap-text style (s="1" on column B in the source
# file); writing .Value only does not disturb that formatting.

$ws.Cells.Item(38,1).Value = 'login_invalid_phone_title'
$ws.Cells.Item(38,2).Value = 'Invalid Phone Number'
$ws.Cells.Item(38,3).Value = 'अमान्य फ़ोन नंबर'
$ws.Cells.Item(38,4).Value = 'చెల్లని ఫోన్ నంబర్'

$ws.Cells.Item(39,1).Value = 'login_invalid_phone_desc'
$ws.Cells.Item(39,2).Value = 'Please enter a valid 10-digit phone number'
$ws.Cells.Item(39,3).Value = 'कृपया मान्य 10 अंकों का फ़ोन नंबर दर्ज करें'
$ws.Cells.Item(39,4).Value = 'దయచేసి సరైన 10 అంకెల ఫోన్ నంబర్ నమోదు చేయండి'

$ws.Cells.Item(40,1).Value = 'login_success_title'
$ws.Cells.Item(40,2).Value = 'Login Successful! 👋'
$ws.Cells.Item(40,3).Value = 'लॉगिन सफल! 👋'
$ws.Cells.Item(40,4).Value = 'లాగిన్ విజయవంతం! 👋'

$ws.Cells.Item(41,1).Value = 'login_success_desc'
$ws.Cells.Item(41,2).Value = 'Welcome back'
$ws.Cells.Item(41,3).Value = 'वापसी पर स्वागत है'
$ws.Cells.Item(41,4).Value = 'మళ్లీ స్వాగతం'

$ws.Cells.Item(42,1).Value = 'login_account_not_found_title'
$ws.Cells.Item(42,2).Value = 'Account Not Found'
$ws.Cells.Item(42,3).Value = 'खाता नहीं मिला'
$ws.Cells.Item(42,4).Value = 'ఖాతా కనబడలేదు'

$ws.Cells.Item(43,1).Value = 'login_account_not_found_desc'
$ws.Cells.Item(43,2).Value = 'This number is not registered. Please sign up first.'
$ws.Cells.Item(43,3).Value = 'यह नंबर पंजीकृत नहीं है। कृपया पहले साइन अप करें।'
$ws.Cells.Item(43,4).Value = 'ఈ నంబర్ నమోదు కాలేదు. దయచేసి ముందుగా నమోదు చేసుకోండి.'

$ws.Cells.Item(44,1).Value = 'login_phone_not_verified_title'
$ws.Cells.Item(44,2).Value = 'Phone Not Verified'
$ws.Cells.Item(44,3).Value = 'फ़ोन सत्यापित नहीं है'
$ws.Cells.Item(44,4).Value = 'ఫోన్ ధృవీకరించబడలేదు'

$ws.Cells.Item(45,1).Value = 'login_phone_not_verified_desc'
$ws.Cells.Item(45,2).Value = 'Please verify your phone number'
$ws.Cells.Item(45,3).Value = 'कृपया अपना फ़ोन नंबर सत्यापित करें'
$ws.Cells.Item(45,4).Value = 'దయచేసి మీ ఫోన్ నంబర్ ధృవీకరించండి'

$ws.Cells.Item(46,1).Value = 'login_failed_title'
$ws.Cells.Item(46,2).Value = 'Login Failed'
$ws.Cells.Item(46,3).Value = 'लॉगिन असफल'
$ws.Cells.Item(46,4).Value = 'లాగిన్ విఫలమైంది'

$ws.Cells.Item(47,1).Value = 'login_network_error_title'
$ws.Cells.Item(47,2).Value = 'Network Error'
$ws.Cells.Item(47,3).Value = 'नेटवर्क त्रुटि'
$ws.Cells.Item(47,4).Value = 'నెట్‌వర్క్ లోపం'

$ws.Cells.Item(48,1).Value = 'login_network_error_desc'
$ws.Cells.Item(48,2).Value = 'Failed to login. Please check your connection.'
$ws.Cells.Item(48,3).Value = 'लॉगिन असफल। कृपया कनेक्शन जांचें.'
$ws.Cells.Item(48,4).Value = 'లాగిన్ చేయలేకపోయాము. దయచేసి కనెక్షన్ తనిఖీ చేయండి.'

$ws.Cells.Item(49,1).Value = 'login_title'
$ws.Cells.Item(49,2).Value = 'Let''s Sign In...'
$ws.Cells.Item(49,3).Value = 'लॉगिन करें...'
$ws.Cells.Item(49,4).Value = 'లాగిన్ అవ్వండి...'

$ws.Cells.Item(50,1).Value = 'login_subtitle'
$ws.Cells.Item(50,2).Value = 'Enter your phone number to access your account'
$ws.Cells.Item(50,3).Value = 'अपने खाते तक पहुंचने के लिए अपना फोन नंबर दर्ज करें'
$ws.Cells.Item(50,4).Value = 'మీ ఖాతాను ఉపయోగించడానికి మీ ఫోన్ నంబర్ నమోదు చేయండి'

$ws.Cells.Item(51,1).Value = 'placeholder_phone'
$ws.Cells.Item(51,2).Value = 'Phone number'
$ws.Cells.Item(51,3).Value = 'फोन नंबर'
$ws.Cells.Item(51,4).Value = 'ఫోన్ నంబర్'

$ws.Cells.Item(52,1).Value = 'login_note_registered_phone'
$ws.Cells.Item(52,2).Value = 'Note: Make sure you registered this phone number before signing in.'
$ws.Cells.Item(52,3).Value = 'नोट: लॉगिन से पहले यह फोन नंबर पंजीकृत होना चाहिए।'
$ws.Cells.Item(52,4).Value = 'గమనిక: లాగిన్ అవ్వడానికి ముందు ఈ ఫోన్ నంబర్ నమోదు చేసి ఉండాలి.'

$ws.Cells.Item(53,1).Value = 'sign_in'
$ws.Cells.Item(53,2).Value = 'Sign In'
$ws.Cells.Item(53,3).Value = 'साइन इन'
$ws.Cells.Item(53,4).Value = 'లాగిన్'

$ws.Cells.Item(54,1).Value = 'register'
$ws.Cells.Item(54,2).Value = 'Register'
$ws.Cells.Item(54,3).Value = 'रजिस्टर'
$ws.Cells.Item(54,4).Value = 'నమోదు చేయండి'

$ws.Cells.Item(55,1).Value = 'already_have_account'
$ws.Cells.Item(55,2).Value = 'Don''t have an account?'
$ws.Cells.Item(55,3).Value = 'खाता नहीं है?'
$ws.Cells.Item(55,4).Value = 'ఖాతా లేదా?'

$ws.Cells.Item(56,1).Value = 'otp_title_enter'
$ws.Cells.Item(56,2).Value = 'Enter Verification code'
$ws.Cells.Item(56,3).Value = 'सत्यापन कोड दर्ज करें'
$ws.Cells.Item(56,4).Value = 'ధృవీకరణ కోడ్ నమోదు చేయండి'

$ws.Cells.Item(57,1).Value = 'otp_subtitle'
$ws.Cells.Item(57,2).Value = 'Please enter 4 digit verification code sent to'
$ws.Cells.Item(57,3).Value = 'कृपया भेजा गया 4 अंकों का सत्यापन कोड दर्ज करें'
$ws.Cells.Item(57,4).Value = 'మీకు పంపబడిన 4 అంకెల ధృవీకరణ కోడ్ నమోదు చేయండి'

$ws.Cells.Item(58,1).Value = 'otp_edit'
$ws.Cells.Item(58,2).Value = 'Edit'
$ws.Cells.Item(58,3).Value = 'संपादित करें'
$ws.Cells.Item(58,4).Value = 'సవరించండి'

# Row 59 (otp_timer_prefix): B/C/D are literal numeric placeholders (0), not text
$ws.Cells.Item(59,1).Value = 'otp_timer_prefix'
$ws.Cells.Item(59,2).Value = 0
$ws.Cells.Item(59,3).Value = 0
$ws.Cells.Item(59,4).Value = 0

$ws.Cells.Item(60,1).Value = 'otp_didnt_receive'
$ws.Cells.Item(60,2).Value = 'Didn''t receive OTP?'
$ws.Cells.Item(60,3).Value = 'OTP नहीं मिला?'
$ws.Cells.Item(60,4).Value = 'OTP అందలేదా?'

$ws.Cells.Item(61,1).Value = 'otp_resend'
$ws.Cells.Item(61,2).Value = 'Resend Code'
$ws.Cells.Item(61,3).Value = 'कोड पुनः भेजें'
$ws.Cells.Item(61,4).Value = 'కోడ్ మళ్లీ పంపండి'

$ws.Cells.Item(62,1).Value = 'otp_submit'
$ws.Cells.Item(62,2).Value = 'Submit'
$ws.Cells.Item(62,3).Value = 'सबमिट करें'
$ws.Cells.Item(62,4).Value = 'సమర్పించండి'

$ws.Cells.Item(63,1).Value = 'alert_invalid_otp_title'
$ws.Cells.Item(63,2).Value = 'Invalid OTP'
$ws.Cells.Item(63,3).Value = 'अमान्य OTP'
$ws.Cells.Item(63,4).Value = 'చెల్లని OTP'

$ws.Cells.Item(64,1).Value = 'alert_invalid_otp_desc'
$ws.Cells.Item(64,2).Value = 'Please enter complete 4-digit OTP'
$ws.Cells.Item(64,3).Value = 'कृपया पूरा 4 अंकों का OTP दर्ज करें'
$ws.Cells.Item(64,4).Value = 'దయచేసి పూర్తి 4 అంకెల OTP నమోదు చేయండి'

$ws.Cells.Item(65,1).Value = 'alert_resend_success_title'
$ws.Cells.Item(65,2).Value = 'Success'
$ws.Cells.Item(65,3).Value = 'सफल'
$ws.Cells.Item(65,4).Value = 'విజయం'

$ws.Cells.Item(66,1).Value = 'alert_resend_success_desc'
$ws.Cells.Item(66,2).Value = 'OTP has been resent to your phone number'
$ws.Cells.Item(66,3).Value = 'OTP आपके फोन नंबर पर फिर से भेजा गया है'
$ws.Cells.Item(66,4).Value = 'OTP మీ ఫోన్ నంబర్‌కు మళ్లీ పంపబడింది'

$ws.Cells.Item(67,1).Value = 'alert_resend_failed_title'
$ws.Cells.Item(67,2).Value = 'Error'
$ws.Cells.Item(67,3).Value = 'त्रुटि'
$ws.Cells.Item(67,4).Value = 'లోపం'

$ws.Cells.Item(68,1).Value = 'alert_resend_failed_desc'
$ws.Cells.Item(68,2).Value = 'Failed to resend OTP. Please try again.'
$ws.Cells.Item(68,3).Value = 'OTP पुनः भेजने में विफल। कृपया पुनः प्रयास करें।'
$ws.Cells.Item(68,4).Value = 'OTP మళ్లీ పంపడంలో విఫలమైంది. దయచేసి మళ్లీ ప్రయత్నించండి.'

$ws.Cells.Item(69,1).Value = 'alert_verification_failed_title'
$ws.Cells.Item(69,2).Value = 'Verification Failed'
$ws.Cells.Item(69,3).Value = 'सत्यापन असफल'
$ws.Cells.Item(69,4).Value = 'ధృవీకరణ విఫలమైంది'

$ws.Cells.Item(70,1).Value = 'alert_verification_failed_desc'
$ws.Cells.Item(70,2).Value = 'Invalid OTP. Please try again.'
$ws.Cells.Item(70,3).Value = 'अमान्य OTP। कृपया पुनः प्रयास करें।'
$ws.Cells.Item(70,4).Value = 'చెల్లని OTP. దయచేసి మళ్లీ ప్రయత్నించండి.'

$ws.Cells.Item(71,1).Value = 'alert_verify_error_title'
$ws.Cells.Item(71,2).Value = 'Error'
$ws.Cells.Item(71,3).Value = 'त्रुटि'
$ws.Cells.Item(71,4).Value = 'లోపం'

$ws.Cells.Item(72,1).Value = 'alert_verify_error_desc'
$ws.Cells.Item(72,2).Value = 'Failed to verify OTP. Please try again.'
$ws.Cells.Item(72,3).Value = 'OTP सत्यापन में विफल। कृपया पुनः प्रयास करें।'
$ws.Cells.Item(72,4).Value = 'OTP ధృవీకరణ విఫలమైంది. దయచేసి మళ్లీ ప్రయత్నించండి.'

$ws.Cells.Item(73,1).Value = 'verify_title_prefix'
$ws.Cells.Item(73,2).Value = 'Enter'
$ws.Cells.Item(73,3).Value = 'दर्ज करें'
$ws.Cells.Item(73,4).Value = 'నమోదు చేయండి'

$ws.Cells.Item(74,1).Value = 'verify_title_main'
$ws.Cells.Item(74,2).Value = 'Verification code'
$ws.Cells.Item(74,3).Value = 'सत्यापन कोड'
$ws.Cells.Item(74,4).Value = 'ధృవీకరణ కోడ్'

$ws.Cells.Item(75,1).Value = 'verify_subtitle_send_otp'
$ws.Cells.Item(75,2).Value = 'We will send an OTP to your registered phone number'
$ws.Cells.Item(75,3).Value = 'हम आपके पंजीकृत फ़ोन नंबर पर OTP भेजेंगे'
$ws.Cells.Item(75,4).Value = 'మీ నమోదు చేసిన ఫోన్ నంబర్‌కు OTP పంపబడుతుంది'

$ws.Cells.Item(76,1).Value = 'verify_label_phone'
$ws.Cells.Item(76,2).Value = 'Phone number'
$ws.Cells.Item(76,3).Value = 'फ़ोन नंबर'
$ws.Cells.Item(76,4).Value = 'ఫోన్ నంబర్'

$ws.Cells.Item(77,1).Value = 'verify_button_send_otp'
$ws.Cells.Item(77,2).Value = 'Send OTP'
$ws.Cells.Item(77,3).Value = 'OTP भेजें'
$ws.Cells.Item(77,4).Value = 'OTP పంపండి'

$ws.Cells.Item(78,1).Value = 'verify_otp_send_failed_title'
$ws.Cells.Item(78,2).Value = 'Error'
$ws.Cells.Item(78,3).Value = 'त्रुटि'
$ws.Cells.Item(78,4).Value = 'లోపం'

$ws.Cells.Item(79,1).Value = 'verify_otp_send_failed_desc'
$ws.Cells.Item(79,2).Value = 'Failed to send OTP. Please try again.'
$ws.Cells.Item(79,3).Value = 'OTP भेजने में विफल। कृपया पुनः प्रयास करें।'
$ws.Cells.Item(79,4).Value = 'OTP పంపడంలో విఫలమైంది. దయచేసి మళ్లీ ప్రయత్నించండి.'

$ws.Cells.Item(80,1).Value = 'verify_network_error_desc'
$ws.Cells.Item(80,2).Value = 'Failed to send OTP. Please check your connection.'
$ws.Cells.Item(80,3).Value = 'OTP भेजने में विफल। कृपया अपना कनेक्शन जांचें।'
$ws.Cells.Item(80,4).Value = 'OTP పంపడంలో విఫలమైంది. దయచేసి మీ కనెక్షన్ తనిఖీ చేయండి.'

# View-state updates: move the active selection the same way the author left it
# (C41 -> C91) and best-effort scroll the window toward the new rows. (Window
# chrome/scroll geometry is session/runtime state that this host does not persist
# back into sheetView/workbookView on save, so these calls are best-effort only.)
$ws.Activate()
try { $excel.ActiveWindow.ScrollRow = 60 } catch { }
try { $excel.ActiveWindow.ScrollColumn = 1 } catch { }
$ws.Range("C91").Select()

